# Append the 2025-10-18 allocation row (A47:C47) to the profit/allocation sheet,
# mirroring the prior daily rows (Date text, BTC fraction, KAS fraction).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces the date-looking string to be stored as literal
# text (matching how the existing Date column cells are stored) instead of
# being auto-converted into a date serial value.
$ws.Range("A47").Value = "'10/18/2025"
$ws.Range("B47").Value = 0.1870601725438363
$ws.Range("C47").Value = 0.8129398274561637
